# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Betarraga, Femacal de La Calera) at the
# top of the data block (rows 543:544), pushing the existing rows down by
# two (old row 543 -> new row 545, ... old row 576 -> new row 578).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything in the data block down by inserting 2 blank rows.
$ws.Rows("543:544").Insert()

# New row 543 - "Primera" grade for the new week (2022-02-18 -> 44610)
$ws.Range("A543").Value = 3
$ws.Range("B543").Value = "Femacal de La Calera"
$ws.Range("C543").Value = "Coquimbo"
$ws.Range("D543").Value = 44610
$ws.Range("E543").Value = 5
$ws.Range("F543").Value = 100114014
$ws.Range("G543").Value = "Betarraga"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 1200
$ws.Range("K543").Value = 600
$ws.Range("L543").Value = 600
$ws.Range("M543").Value = 600
$ws.Range("N543").Value = "$/paquete 4 unidades"
$ws.Range("O543").Value = "Provincia de Quillota"
$ws.Range("P543").Value = 150
$ws.Range("Q543").Value = 4
$ws.Range("R543").Value = "Hortaliza"

# New row 544 - "Segunda" grade for the new week (2022-02-18 -> 44610)
$ws.Range("A544").Value = 3
$ws.Range("B544").Value = "Femacal de La Calera"
$ws.Range("C544").Value = "Coquimbo"
$ws.Range("D544").Value = 44610
$ws.Range("E544").Value = 5
$ws.Range("F544").Value = 100114014
$ws.Range("G544").Value = "Betarraga"
$ws.Range("H544").Value = "Sin especificar"
$ws.Range("I544").Value = "Segunda"
$ws.Range("J544").Value = 850
$ws.Range("K544").Value = 400
$ws.Range("L544").Value = 400
$ws.Range("M544").Value = 400
$ws.Range("N544").Value = "$/paquete 4 unidades"
$ws.Range("O544").Value = "Provincia de Quillota"
$ws.Range("P544").Value = 100
$ws.Range("Q544").Value = 4
$ws.Range("R544").Value = "Hortaliza"
